$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ROW2
$ws.Range("F2").Value = 24.50000000000039
$ws.Range("H2").Value = 0.0007320602578930746
$ws.Range("I2").Value = 0.0007320602578930746
$ws.Range("L2").Value = 36.28634598962337
$ws.Range("M2").Value = "[15.378274804880917, 57.194417174365825]"
$ws.Range("N2").Value = 0.001075576388872346
$ws.Range("O2").Value = 0.001075576388872346
$ws.Range("P2").Value = 1.628973968528041
$ws.Range("Q2").Value = "[0.8868159442565782, 2.371131992799504]"
$ws.Range("R2").Value = 0.00006148203804934305
$ws.Range("S2").Value = 0.00006148203804934305
$ws.Range("T2").Value = 56.13252117629771
$ws.Range("U2").Value = "[43.49158257218285, 68.77345978041257]"
$ws.Range("V2").Value = 0.0000000000152193813107715
$ws.Range("W2").Value = 0.0000000000152193813107715
$ws.Range("X2").Value = 18.14814814814844
$ws.Range("Y2").Value = 15.25425425425449
$ws.Range("Z2").Value = 21.04204204204238

# ROW3
$ws.Range("F3").Value = 24.50000000000039
$ws.Range("H3").Value = 0.0000000612108659314714
$ws.Range("I3").Value = 0.0000000612108659314714
$ws.Range("L3").Value = 54.60456002831246
$ws.Range("M3").Value = "[37.40034811980071, 71.8087719368242]"
$ws.Range("N3").Value = 0.00000008145305718798568
$ws.Range("O3").Value = 0.00000008145305718798568
$ws.Range("P3").Value = 1.50318447288881
$ws.Range("Q3").Value = "[1.1258159859711174, 1.8805529598065025]"
$ws.Range("R3").Value = 0.0000000003193108000232314
$ws.Range("S3").Value = 0.0000000003193108000232314
$ws.Range("T3").Value = 63.27904756364193
$ws.Range("U3").Value = "[52.08502833110293, 74.47306679618093]"
$ws.Range("V3").Value = 0.000000000000007771561172376096
$ws.Range("W3").Value = 0.000000000000007771561172376096
$ws.Range("X3").Value = 18.63863863863894
$ws.Range("Y3").Value = 17.16716716716745
$ws.Range("Z3").Value = 20.11011011011043

# ROW4
$ws.Range("F4").Value = 24.50000000000039
$ws.Range("H4").Value = 0.000000002455148417901398
$ws.Range("I4").Value = 0.000000002455148417901398
$ws.Range("L4").Value = 55.31621884439546
$ws.Range("M4").Value = "[39.690169166448314, 70.9422685223426]"
$ws.Range("N4").Value = 0.00000000654124865384631
$ws.Range("O4").Value = 0.00000000654124865384631
$ws.Range("P4").Value = 1.616395018964118
$ws.Range("Q4").Value = "[1.2893423303021176, 1.9434477076261185]"
$ws.Range("R4").Value = 0.0000000000006001865671123596
$ws.Range("S4").Value = 0.0000000000006001865671123596
$ws.Range("T4").Value = 60.94714786754028
$ws.Range("U4").Value = "[50.96754837878679, 70.92674735629377]"
$ws.Range("V4").Value = 0.0000000000000006661338147750939
$ws.Range("W4").Value = 0.0000000000000006661338147750939
$ws.Range("X4").Value = 18.19719719719749
$ws.Range("Y4").Value = 16.92192192192219
$ws.Range("Z4").Value = 19.47247247247278

# ROW5
$ws.Range("F5").Value = 24.50000000000039
$ws.Range("H5").Value = 0.000000003972279838393433
$ws.Range("I5").Value = 0.000000003972279838393433
$ws.Range("L5").Value = 52.96267771940813
$ws.Range("M5").Value = "[37.47654543188729, 68.44881000692897]"
$ws.Range("N5").Value = 0.00000001492713996853468
$ws.Range("O5").Value = 0.00000001492713996853468
$ws.Range("P5").Value = 1.452868674633118
$ws.Range("Q5").Value = "[1.1132370364071944, 1.7925003128590413]"
$ws.Range("R5").Value = 0.00000000004453504232060368
$ws.Range("S5").Value = 0.00000000004453504232060368
$ws.Range("T5").Value = 55.29268275347487
$ws.Range("U5").Value = "[45.5866713998947, 64.99869410705504]"
$ws.Range("V5").Value = 0.000000000000005773159728050814
$ws.Range("W5").Value = 0.000000000000005773159728050814
$ws.Range("X5").Value = 18.83483483483514
$ws.Range("Y5").Value = 17.51051051051079
$ws.Range("Z5").Value = 20.15915915915948

# ROW6
$ws.Range("F6").Value = 24.50000000000039
$ws.Range("H6").Value = 0.0000001748299912618378
$ws.Range("I6").Value = 0.0000001748299912618378
$ws.Range("L6").Value = 56.26716099496016
$ws.Range("M6").Value = "[34.190935413252106, 78.34338657666822]"
$ws.Range("N6").Value = 0.00000589078212143157
$ws.Range("O6").Value = 0.00000589078212143157
$ws.Range("P6").Value = 1.239026532046425
$ws.Range("Q6").Value = "[0.81134224687304, 1.6667108172198093]"
$ws.Range("R6").Value = 0.0000005485383747227957
$ws.Range("S6").Value = 0.0000005485383747227957
$ws.Range("T6").Value = 52.18739570201232
$ws.Range("U6").Value = "[40.15278811248498, 64.22200329153966]"
$ws.Range("V6").Value = 0.00000000003020295125111261
$ws.Range("W6").Value = 0.00000000003020295125111261
$ws.Range("X6").Value = 19.66866866866899
$ws.Range("Y6").Value = 18.00100100100129
$ws.Range("Z6").Value = 21.33633633633668

# ROW7
$ws.Range("F7").Value = 24.50000000000039
$ws.Range("H7").Value = 0.00000001713876696118888
$ws.Range("I7").Value = 0.00000001713876696118888
$ws.Range("L7").Value = 55.32753068247587
$ws.Range("M7").Value = "[34.938125206496096, 75.71693615845564]"
$ws.Range("N7").Value = 0.000001926675075036854
$ws.Range("O7").Value = 0.000001926675075036854
$ws.Range("P7").Value = 1.037763339023655
$ws.Range("Q7").Value = "[0.6478159025420398, 1.427710775505271]"
$ws.Range("R7").Value = 0.000002749665964252301
$ws.Range("S7").Value = 0.000002749665964252301
$ws.Range("T7").Value = 61.93232661457178
$ws.Range("U7").Value = "[51.24928693922419, 72.61536628991936]"
$ws.Range("V7").Value = 0.00000000000000333066907387547
$ws.Range("W7").Value = 0.00000000000000333066907387547
$ws.Range("X7").Value = 20.45345345345378
$ws.Range("Y7").Value = 18.93293293293323
$ws.Range("Z7").Value = 21.97397397397432

# ROW8
$ws.Range("F8").Value = 23.75000000000027
$ws.Range("H8").Value = 0.000004234819360515729
$ws.Range("I8").Value = 0.000004234819360515729
$ws.Range("L8").Value = 45.11354869580541
$ws.Range("M8").Value = "[23.914210064589643, 66.31288732702117]"
$ws.Range("N8").Value = 0.00009466967610483046
$ws.Range("O8").Value = 0.00009466967610483046
$ws.Range("P8").Value = 0.9371317425122712
$ws.Range("Q8").Value = "[0.45913165908319353, 1.4151318259413488]"
$ws.Range("R8").Value = 0.0002735309226040705
$ws.Range("S8").Value = 0.0002735309226040705
$ws.Range("T8").Value = 51.35637203673472
$ws.Range("U8").Value = "[40.18537418495756, 62.52736988851188]"
$ws.Range("V8").Value = 0.000000000005471845199167547
$ws.Range("W8").Value = 0.000000000005471845199167547
$ws.Range("X8").Value = 20.20770770770794
$ws.Range("Y8").Value = 18.40090090090111
$ws.Range("Z8").Value = 22.01451451451477

# ROW9
$ws.Range("F9").Value = 23.75000000000027
$ws.Range("H9").Value = 0.0000000001439515173728978
$ws.Range("I9").Value = 0.0000000001439515173728978
$ws.Range("L9").Value = 65.15585881389696
$ws.Range("M9").Value = "[44.58441219854194, 85.72730542925197]"
$ws.Range("N9").Value = 0.00000008525408246029542
$ws.Range("O9").Value = 0.00000008525408246029542
$ws.Range("P9").Value = 0.6855527512338089
$ws.Range("Q9").Value = "[0.3836579616996536, 0.9874475407679641]"
$ws.Range("R9").Value = 0.00003747015066180026
$ws.Range("S9").Value = 0.00003747015066180026
$ws.Range("T9").Value = 58.04184725054607
$ws.Range("U9").Value = "[47.534556570377134, 68.549137930715]"
$ws.Range("V9").Value = 0.00000000000001643130076445232
$ws.Range("W9").Value = 0.00000000000001643130076445232
$ws.Range("X9").Value = 21.1586586586589
$ws.Range("Y9").Value = 20.01751751751774
$ws.Range("Z9").Value = 22.29979979980006

# ROW10
$ws.Range("F10").Value = 23.75000000000027
$ws.Range("H10").Value = 0.00000001227425872052379
$ws.Range("I10").Value = 0.00000001227425872052379
$ws.Range("L10").Value = 58.32496386598785
$ws.Range("M10").Value = "[37.352566032782576, 79.29736169919313]"
$ws.Range("N10").Value = 0.00000121521093010557
$ws.Range("O10").Value = 0.00000121521093010557
$ws.Range("P10").Value = 1.062921238151502
$ws.Range("Q10").Value = "[0.6855527512338089, 1.4402897250691957]"
$ws.Range("R10").Value = 0.000000952196767434188
$ws.Range("S10").Value = 0.000000952196767434188
$ws.Range("T10").Value = 56.52103367497322
$ws.Range("U10").Value = "[45.46133004007094, 67.5807373098755]"
$ws.Range("V10").Value = 0.0000000000002087219286295294
$ws.Range("W10").Value = 0.0000000000002087219286295294
$ws.Range("X10").Value = 19.73223223223246
$ws.Range("Y10").Value = 18.30580580580601
$ws.Range("Z10").Value = 21.1586586586589

# ROW11
$ws.Range("F11").Value = 23.75000000000027
$ws.Range("H11").Value = 0.000002905614715498217
$ws.Range("I11").Value = 0.000002905614715498217
$ws.Range("L11").Value = 45.17917839763882
$ws.Range("M11").Value = "[24.066514715096147, 66.2918420801815]"
$ws.Range("N11").Value = 0.00008772713535587506
$ws.Range("O11").Value = 0.00008772713535587506
$ws.Range("P11").Value = 0.9497106920761942
$ws.Range("Q11").Value = "[0.4842895582110396, 1.4151318259413488]"
$ws.Range("R11").Value = 0.000165441554323742
$ws.Range("S11").Value = 0.000165441554323742
$ws.Range("T11").Value = 55.9201339972704
$ws.Range("U11").Value = "[44.96222176414899, 66.87804623039182]"
$ws.Range("V11").Value = 0.0000000000002184918912462308
$ws.Range("W11").Value = 0.0000000000002184918912462308
$ws.Range("X11").Value = 20.16016016016039
$ws.Range("Y11").Value = 18.40090090090111
$ws.Range("Z11").Value = 21.91941941941967

# ROW12
$ws.Range("F12").Value = 23.75000000000027
$ws.Range("H12").Value = 0.00001219201659430347
$ws.Range("I12").Value = 0.00001219201659430347
$ws.Range("L12").Value = 44.26029098757413
$ws.Range("M12").Value = "[22.045769203637292, 66.47481277151097]"
$ws.Range("N12").Value = 0.0002240820566856705
$ws.Range("O12").Value = 0.0002240820566856705
$ws.Range("P12").Value = 1.037763339023655
$ws.Range("Q12").Value = "[0.5220264069028087, 1.553500271144502]"
$ws.Range("R12").Value = 0.0001978673728328939
$ws.Range("S12").Value = 0.0001978673728328939
$ws.Range("T12").Value = 53.0236194444773
$ws.Range("U12").Value = "[41.430388804538396, 64.6168500844162]"
$ws.Range("V12").Value = 0.000000000006379119454891224
$ws.Range("W12").Value = 0.000000000006379119454891224
$ws.Range("X12").Value = 19.82732732732756
$ws.Range("Y12").Value = 17.87787787787808
$ws.Range("Z12").Value = 21.77677677677703
